$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
  @("a", 1042.0),
  @("a", 1042.0),
  @("a", 1042.0),
  @("a", 1042.0),
  @("a", 521.0),
  @("l", 0.0),
  @("l", 1042.0),
  @("a", 1042.0),
  @("k", 1563.0),
  @("a", 1042.0),
  @("a", 1042.0),
  @("l", 1042.0),
  @("k", 1042.0),
  @("k", 0.0),
  @("k", 4689.0),
  @("k", 0.0),
  @("g", 2084.0),
  @("a", 1042.0),
  @("l", 1042.0),
  @("a", 1042.0),
  @("k", 1042.0),
  @("a", 0.0),
  @("l", 1042.0),
  @("k", 1042.0),
  @("l", 1563.0),
  @("a", 1042.0),
  @("a", 1042.0),
  @("a", 5210.0),
  @("a", 0.0),
  @("l", 0.0),
  @("l", 0.0),
  @("l", 3647.0),
  @("l", 2084.0),
  @("a", 1042.0),
  @("l", 0.0),
  @("l", 1042.0),
  @("l", 1042.0),
  @("a", 1563.0),
  @("l", 1042.0),
  @("l", 1042.0),
  @("a", 1042.0),
  @("a", 0.0),
  @("l", 0.0),
  @("l", 0.0),
  @("l", 0.0),
  @("a", 0.0),
  @("a", 2605.0),
  @("a", 0.0),
  @("a", 0.0),
  @("a", 1042.0),
  @("a", 3126.0),
  @("l", 1563.0),
  @("l", 2605.0)
)

$startRow = 34
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
